$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values regenerated to filter save games (rows 2-8, columns B-G)
$data = @{
    2 = @{ B = 0.6545652718822623;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; F = 0; G = 6.038307959104277 }
    3 = @{ B = 3.272327238179451;   C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; F = 0; G = 5.582307763322248 }
    4 = @{ B = 3.272327238179451;   C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; F = 0; G = 6.15379541431027 }
    5 = @{ B = 1.445647641019636;   C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; F = 1; G = 4.327115817150455 }
    6 = @{ B = 0.01253208636536152; C = 0.04103571897497393; D = 0.1496068669990043; E = 0.5333859586016987; F = 1; G = 0.7365606309410384 }
    7 = @{ B = 3.272327238179451;   C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; F = 1; G = 5.582307763322248 }
    8 = @{ B = 0.2881169905109251;  C = 0.3048912486333797; D = 0.7210945179870265; E = 0.5333859586016987; F = 1; G = 1.84748871573303 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
